# Update ticket-price column (G) from numeric (cents) to descriptive text,
# and bump a handful of "want to go" counters (F) to match refreshed scrape data.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "不可售"
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "已售罄"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 6).Value = 1741
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "55"
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "68"
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "60"
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "80"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(8, 6).Value = 711
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "48"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "65"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "60"
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "60"
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "60"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "68"
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(14, 6).Value = 1244
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "35"
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "60"
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "60"
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "63"
$ws.Cells.Item(17, 7).Style = "Normal"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "60"
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "58"
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "58"
$ws.Cells.Item(20, 7).Style = "Normal"
$ws.Cells.Item(21, 6).Value = 178
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "55"
$ws.Cells.Item(21, 7).Style = "Normal"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "55"
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "48"
$ws.Cells.Item(23, 7).Style = "Normal"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "60"
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "40"
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "60"
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "39"
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "60"
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "60"
$ws.Cells.Item(29, 7).Style = "Normal"
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "69"
$ws.Cells.Item(30, 7).Style = "Normal"

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "不可售"
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "180"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "380"
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "100"
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "50"
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(7, 6).Value = 818
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "380"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "420"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "380"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "280"
$ws.Cells.Item(10, 7).Style = "Normal"

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "不可售"
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "不可售"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "180"
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "已售罄"
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(6, 6).Value = 1741
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "55"
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "380"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "68"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "60"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "80"
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(11, 6).Value = 711
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "48"
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "65"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "100"
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "60"
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "60"
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "60"
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "68"
$ws.Cells.Item(17, 7).Style = "Normal"
$ws.Cells.Item(18, 6).Value = 1244
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "35"
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "50"
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "60"
$ws.Cells.Item(20, 7).Style = "Normal"
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "60"
$ws.Cells.Item(21, 7).Style = "Normal"
$ws.Cells.Item(22, 6).Value = 818
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "380"
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "63"
$ws.Cells.Item(23, 7).Style = "Normal"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "60"
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "58"
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "58"
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(27, 6).Value = 178
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "55"
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "55"
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "48"
$ws.Cells.Item(29, 7).Style = "Normal"
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "60"
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "420"
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "40"
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "60"
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "39"
$ws.Cells.Item(34, 7).Style = "Normal"
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "60"
$ws.Cells.Item(35, 7).Style = "Normal"
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "380"
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "60"
$ws.Cells.Item(37, 7).Style = "Normal"
$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "69"
$ws.Cells.Item(38, 7).Style = "Normal"
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "280"
$ws.Cells.Item(39, 7).Style = "Normal"
